$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct verse references: F11's border style changes (s=1 -> s=3) ---
# since row 12 below it now also holds a letter, so F11 is no longer the
# bottom of its own isolated box. Copy the border formatting from another
# cell that already uses the "open-bottom" style (e.g. L6) so the existing
# style index is reused rather than a new one being created.
$ws.Range("L6").Copy()
$ws.Range("F11").PasteSpecial(-4122)

# --- Added an additional word: "YOUTH" across B12:F12 (F12 already has H) ---
$newLetters = @{ "B12" = "Y"; "C12" = "O"; "D12" = "U"; "E12" = "T" }
foreach ($addr in $newLetters.Keys) {
    # Copy format from an existing fully-bordered answer cell (style s=1)
    # so the new cells reuse that style rather than minting a new one.
    $ws.Range("B3").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $newLetters[$addr]
}

# --- Update the saved selection / active cell ---
$ws.Range("L10").Select()
